# إضافة حدث جديد في Card20 by admin at 2025-12-16 11:06:47
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Duplicate the last log row (34) down into the new row (35) so the new
# entry inherits the same "card/date/event/correction/serviced by" shape.
$ws.Range("A34:O34").Copy($ws.Range("A35:O35"))

# Back-fill the now-historical row 34's empty measurement columns (B:K)
# with the "nan" placeholder used elsewhere in the sheet.
$ws.Range("B34:K34").Value = "nan"
